# Refactored Ask Ana translation management
# Adds a new "lang_config" worksheet (at the end of the tab list) that
# tracks which languages the blog automation should translate into, and
# which translation stages ("post" / "run") are enabled per language.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end
# of the tab strip (mirrors blocked_domains -> lang_config ordering).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$langConfig = $wb.Worksheets.Add($null, $lastSheet)
$langConfig.Name = "lang_config"

# Header row.
$langConfig.Range("A1").Value = "lang"
$langConfig.Range("B1").Value = "post"
$langConfig.Range("C1").Value = "run"
$langConfig.Range("D1").Value = "code"

# Column A first (language display names), matching the shared-string
# insertion order of the source edit.
$langConfig.Range("A2").Value = "German"
$langConfig.Range("A3").Value = "Russian"
$langConfig.Range("A4").Value = "Chinese (simplified)"
$langConfig.Range("A5").Value = "Hindi"
$langConfig.Range("A6").Value = "Spanish"
$langConfig.Range("A7").Value = "French"

# Then the boolean toggle columns.
$langConfig.Range("B2").Value = $true
$langConfig.Range("C2").Value = $true

$langConfig.Range("B3").Value = $true
$langConfig.Range("C3").Value = $true

$langConfig.Range("B4").Value = $true
$langConfig.Range("C4").Value = $true

$langConfig.Range("B5").Value = $false
$langConfig.Range("C5").Value = $false

$langConfig.Range("B6").Value = $false
$langConfig.Range("C6").Value = $false

$langConfig.Range("B7").Value = $true
$langConfig.Range("C7").Value = $false

# Finally the language codes (column D).
$langConfig.Range("D2").Value = "de"
$langConfig.Range("D3").Value = "ru"
$langConfig.Range("D4").Value = "zh"
$langConfig.Range("D5").Value = "hi"
$langConfig.Range("D6").Value = "es"
$langConfig.Range("D7").Value = "fr"

# Match the saved selection/active-cell on the new sheet.
$null = $langConfig.Range("D8").Select()
